{"js": "// Fill in the \"Actual runtime:\" / \"Actual ranking:\" measurements that were\n// recorded after the lab was run, for each of the six conversion sections.\n// Each section is identified by its Heading 1 title; the runtime line\n// immediately follows the heading, and the ranking line immediately\n// follows the runtime line.\n\nconst sections = [\n  { heading: \"Decimal to Binary\", runtimeSuffix: \" 1.2311699 milliseconds\", rankingSuffix: \" 1\" },\n  { heading: \"Binary to Decimal\", runtimeSuffix: \"3.213598 milliseconds\", rankingSuffix: \" 4\" },\n  { heading: \"Hexadecimal to Decimal\", runtimeSuffix: \"2.9421 milliseconds\", rankingSuffix: \" 3\" },\n  { heading: \"Decimal to Hexadecimal\", runtimeSuffix: \"10.13061 milliseconds\", rankingSuffix: \"  5\" },\n  { heading: \"Binary to Hexadecimal\", runtimeSuffix: \"11.47149 milliseconds \", rankingSuffix: \" 6\" },\n  { heading: \"Hexadecimal to Binary\", runtimeSuffix: \"1.34908 milliseconds\", rankingSuffix: \" 2\" },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfor (const section of sections) {\n  // Find the Heading 1 paragraph with the matching title.\n  let headingIndex = -1;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].style === \"Heading 1\" && items[i].text.trim() === section.heading) {\n      headingIndex = i;\n      break;\n    }\n  }\n  if (headingIndex === -1) {\n    continue;\n  }\n\n  const runtimePara = items[headingIndex + 1];\n  const rankingPara = items[headingIndex + 2];\n\n  if (runtimePara && runtimePara.text.indexOf(\"Actual runtime:\") !== -1) {\n    runtimePara.insertText(section.runtimeSuffix, \"End\");\n  }\n  if (rankingPara && rankingPara.text.indexOf(\"Actual ranking:\") !== -1) {\n    rankingPara.insertText(section.rankingSuffix, \"End\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fill in the \"Actual runtime:\" / \"Actual ranking:\" measurements that were\n# recorded after the lab was run, for each of the six conversion sections.\n# Each section is identified by its Heading 1 title; the runtime line\n# immediately follows the heading, and the ranking line immediately\n# follows the runtime line.\n\n$d = $word.ActiveDocument\n\n$sections = @(\n    @{ Heading = \"Decimal to Binary\";      RuntimeSuffix = \" 1.2311699 milliseconds\"; RankingSuffix = \" 1\" },\n    @{ Heading = \"Binary to Decimal\";      RuntimeSuffix = \"3.213598 milliseconds\";   RankingSuffix = \" 4\" },\n    @{ Heading = \"Hexadecimal to Decimal\"; RuntimeSuffix = \"2.9421 milliseconds\";     RankingSuffix = \" 3\" },\n    @{ Heading = \"Decimal to Hexadecimal\"; RuntimeSuffix = \"10.13061 milliseconds\";   RankingSuffix = \"  5\" },\n    @{ Heading = \"Binary to Hexadecimal\";  RuntimeSuffix = \"11.47149 milliseconds \";  RankingSuffix = \" 6\" },\n    @{ Heading = \"Hexadecimal to Binary\";  RuntimeSuffix = \"1.34908 milliseconds\";    RankingSuffix = \" 2\" }\n)\n\n$count = $d.Paragraphs.Count\n\nforeach ($section in $sections) {\n    $headingIndex = -1\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $text = $p.Range.Text.TrimEnd([char]13)\n        if (($p.Style.NameLocal -eq \"Heading 1\") -and ($text.Trim() -eq $section.Heading)) {\n            $headingIndex = $i\n            break\n        }\n    }\n\n    if ($headingIndex -eq -1) {\n        continue\n    }\n\n    $runtimePara = $d.Paragraphs.Item($headingIndex + 1)\n    $rankingPara = $d.Paragraphs.Item($headingIndex + 2)\n\n    $runtimeText = $runtimePara.Range.Text.TrimEnd([char]13)\n    if ($runtimeText.Contains(\"Actual runtime:\")) {\n        $runtimePara.Range.InsertAfter($section.RuntimeSuffix)\n    }\n\n    $rankingText = $rankingPara.Range.Text.TrimEnd([char]13)\n    if ($rankingText.Contains(\"Actual ranking:\")) {\n        $rankingPara.Range.InsertAfter($section.RankingSuffix)\n    }\n}\n"}
